$wb = $excel.ActiveWorkbook

$sheetNames = @("Items", "Items - Formatted")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A2").Value  = "green onion Pancakes ☐ (1)"
    $ws.Range("A3").Value  = "Pan Fried Leek Dumplings #T (2)"
    $ws.Range("A4").Value  = "Pork Xiao Long Bao(10) À#122E(10)"
    $ws.Range("A5").Value  = "Q-BAO (5) #E,EL (5)"
    $ws.Range("A6").Value  = "Chicken potstickers 3È#45(6)"
    $ws.Range("A7").Value  = "Tomato Mushroom Steamed dumpli ptkINtA0 (6)"
    $ws.Range("A8").Value  = "Zucchini shrimp dumplings #/2"
    $ws.Range("A9").Value  = "beef stew nodle soup (Non Spicy U4P11(7#)"
    $ws.Range("A10").Value = "dandan noodle ##iE"
    $ws.Range("A11").Value = "banana naan bread ZATRA"
    $ws.Range("A12").Value = "house made plum juice G`$MSH+"
}
